# Clarify "part-time job" instead of "job" in the Project Introduction
# paragraph: "university–job matching platform" -> "university– part time job matching platform"

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "university" + [char]0x2013 + "job"
$find.Replacement.Text = "university" + [char]0x2013 + " part time job"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
